# Merge the three runs that together spell out "<id>p025v_1</id>" into a
# single run (keeping the formatting -- Courier New / color 7f6000 -- of
# the first run, which already carries "<id>") so the final text of that
# run becomes "<id>p025v_1</id>".
#
# Strategy: find the paragraph whose text contains the full "<id>...</id>"
# tag, then use Range surgery (Delete + InsertAfter) rather than a plain
# Find/Replace or Range.Text assignment. Find/Replace & Range.Text both
# correctly merge the runs, but InsertAfter is what actually reproduces
# the exact xml:space="preserve" serialization on the resulting <w:t>.

$d = $word.ActiveDocument

foreach ($para in $d.Paragraphs) {
    $full = $para.Range
    $text = $full.Text

    if ($text.Contains("<id>") -and $text.Contains("</id>")) {
        $openTag = "<id>"

        $start = $full.Start
        $openEnd = $start + $openTag.Length

        # Everything in the paragraph after "<id>", excluding the final
        # paragraph mark (full.End - 1).
        $tailRange = $d.Range($openEnd, $full.End - 1)
        $tailText = $tailRange.Text

        # Remove "p025v_1</id>" (currently split across two more runs)
        # and re-insert it right after "<id>" so it becomes part of that
        # first run.
        $tailRange.Delete()

        $insertPoint = $d.Range($openEnd, $openEnd)
        $insertPoint.InsertAfter($tailText)

        break
    }
}
